# Append two new rows of data (rows 4 and 5) to the GILD bag sheet,
# matching the columns: Date, totalScore, posWordPercentage, negWordPercentage,
# posPhrasePercentage, negPhrasePercentage, ElapsedMs, wordCount, sentenceCount,
# posWordCount, negWordCount, positivePhraseCount, negativePhraseCount, Method

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row4 = @(42606.56040509259, -28, 60, 39, 20, 80, 7710, 11537, 1363, 131, 85, 10, 40, "Bag")
$row5 = @(42606.572256944448, -30, 64, 32, 14, 85, 8801, 13212, 1516, 184, 94, 9, 53, "Bag")

for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(4, $col).Value = $row4[$col - 1]
    $ws.Cells.Item(5, $col).Value = $row5[$col - 1]
}

# Match the date number format used by the existing Date column (A2:A3)
# by copying the existing cell's format rather than assigning a NumberFormat
# string (which would create a brand new custom number format entry).
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
